# Auto-generated edit script: update Leve profit calculation columns (H-N)
# across multiple worksheets, per the scheduled price-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 447.375
$ws.Range("I2").Value = 225.57143
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 225.57143
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -112.57143
$ws.Range("N2").Value = -2226

$ws.Range("H18").Value = 2814.4285
$ws.Range("I18").Value = 1616.8334
$ws.Range("K18").Value = 1616.8334
$ws.Range("M18").Value = -1332.8334

$ws.Range("H70").Value = 89306.336
$ws.Range("I70").Value = 1847.5
$ws.Range("J70").Value = 106798.1
$ws.Range("K70").Value = 5542.5
$ws.Range("L70").Value = 320394.3
$ws.Range("M70").Value = -5272.5
$ws.Range("N70").Value = -320934.3

$ws.Range("H73").Value = 89306.336
$ws.Range("I73").Value = 1847.5
$ws.Range("J73").Value = 106798.1
$ws.Range("K73").Value = 5542.5
$ws.Range("L73").Value = 320394.3
$ws.Range("M73").Value = -4606.5
$ws.Range("N73").Value = -322266.3

$ws.Range("H101").Value = 1118
$ws.Range("I101").Value = 200
$ws.Range("K101").Value = 600
$ws.Range("M101").Value = 1022

$ws.Range("H106").Value = 994.4
$ws.Range("I106").Value = 994.4
$ws.Range("K106").Value = 994.4
$ws.Range("M106").Value = -363.4

$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 3000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1080
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9447.562
$ws.Range("I32").Value = 4282.7886
$ws.Range("J32").Value = 22236.523
$ws.Range("K32").Value = 4282.7886
$ws.Range("L32").Value = 22236.523
$ws.Range("M32").Value = -3995.7886
$ws.Range("N32").Value = -22810.523

$ws.Range("H45").Value = 5580.6665
$ws.Range("I45").Value = 6118.524
$ws.Range("K45").Value = 6118.524
$ws.Range("M45").Value = -5741.524

$ws.Range("H88").Value = 1474.25
$ws.Range("I88").Value = 2050
$ws.Range("J88").Value = 1282.3334
$ws.Range("K88").Value = 2050
$ws.Range("L88").Value = 1282.3334
$ws.Range("M88").Value = -1644
$ws.Range("N88").Value = -2094.3334

$ws.Range("H91").Value = 1474.25
$ws.Range("I91").Value = 2050
$ws.Range("J91").Value = 1282.3334
$ws.Range("K91").Value = 2050
$ws.Range("L91").Value = 1282.3334
$ws.Range("M91").Value = -646
$ws.Range("N91").Value = -4090.3334

$ws.Range("H122").Value = 1302.5238
$ws.Range("I122").Value = 1175.7368
$ws.Range("K122").Value = 3527.2104
$ws.Range("M122").Value = -1077.2104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4978.5386
$ws.Range("I20").Value = 4376.625
$ws.Range("J20").Value = 5941.6
$ws.Range("K20").Value = 4376.625
$ws.Range("L20").Value = 5941.6
$ws.Range("M20").Value = -4129.625
$ws.Range("N20").Value = -6435.6

$ws.Range("H86").Value = 2327.8845
$ws.Range("I86").Value = 1617.7778
$ws.Range("J86").Value = 3925.625
$ws.Range("K86").Value = 1617.7778
$ws.Range("L86").Value = 3925.625
$ws.Range("M86").Value = -494.7778000000001
$ws.Range("N86").Value = -6171.625

$ws.Range("H89").Value = 2327.8845
$ws.Range("I89").Value = 1617.7778
$ws.Range("J89").Value = 3925.625
$ws.Range("K89").Value = 8088.889
$ws.Range("L89").Value = 19628.125
$ws.Range("M89").Value = -2472.889
$ws.Range("N89").Value = -30860.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 46280.332
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H58").Value = 3904.3845
$ws.Range("I58").Value = 3659.28
$ws.Range("J58").Value = 4342.0713
$ws.Range("K58").Value = 3659.28
$ws.Range("L58").Value = 4342.0713
$ws.Range("M58").Value = -3456.28
$ws.Range("N58").Value = -4748.0713

$ws.Range("H105").Value = 886.61536
$ws.Range("I105").Value = 889.6667
$ws.Range("K105").Value = 889.6667
$ws.Range("M105").Value = 857.3333

$ws.Range("H132").Value = 3896.6667
$ws.Range("I132").Value = 3788.7646
$ws.Range("K132").Value = 11366.2938
$ws.Range("M132").Value = -8836.293799999999

$ws.Range("H136").Value = 3904.3845
$ws.Range("I136").Value = 3659.28
$ws.Range("J136").Value = 4342.0713
$ws.Range("K136").Value = 10977.84
$ws.Range("L136").Value = 13026.2139
$ws.Range("M136").Value = -8427.84
$ws.Range("N136").Value = -18126.2139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 22204.64
$ws.Range("I131").Value = 84577.086
$ws.Range("J131").Value = 2508.0789
$ws.Range("K131").Value = 253731.258
$ws.Range("L131").Value = 7524.236699999999
$ws.Range("M131").Value = -248691.258
$ws.Range("N131").Value = -17604.2367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 6666.6665
$ws.Range("I57").Value = 6666.6665
$ws.Range("K57").Value = 6666.6665
$ws.Range("M57").Value = -5846.6665

$ws.Range("H70").Value = 11338
$ws.Range("I70").Value = 4007
$ws.Range("K70").Value = 4007
$ws.Range("M70").Value = -3737

$ws.Range("H73").Value = 11338
$ws.Range("I73").Value = 4007
$ws.Range("K73").Value = 4007
$ws.Range("M73").Value = -3071

$ws.Range("H102").Value = 35715092
$ws.Range("I102").Value = 836.6667
$ws.Range("K102").Value = 836.6667
$ws.Range("M102").Value = 785.3333

$ws.Range("H113").Value = 2988.0476
$ws.Range("I113").Value = 3033.35
$ws.Range("J113").Value = 2082
$ws.Range("K113").Value = 3033.35
$ws.Range("L113").Value = 2082
$ws.Range("M113").Value = -863.3499999999999
$ws.Range("N113").Value = -6422

$ws.Range("H122").Value = 3976.85
$ws.Range("I122").Value = 3632.4614
$ws.Range("J122").Value = 4616.4287
$ws.Range("K122").Value = 10897.3842
$ws.Range("L122").Value = 13849.2861
$ws.Range("M122").Value = -8447.3842
$ws.Range("N122").Value = -18749.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 7000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H104").Value = 38495.5
$ws.Range("J104").Value = 38495.5
$ws.Range("L104").Value = 38495.5
$ws.Range("N104").Value = -45483.5

$ws.Range("H122").Value = 1801.1818
$ws.Range("I122").Value = 1773.3889
$ws.Range("K122").Value = 5320.1667
$ws.Range("M122").Value = -2870.1667

$ws.Range("H132").Value = 2925.111
$ws.Range("I132").Value = 2920.7058
$ws.Range("K132").Value = 8762.117400000001
$ws.Range("M132").Value = -6232.117400000001

$ws.Range("H135").Value = 117250
$ws.Range("J135").Value = 117250
$ws.Range("L135").Value = 117250
$ws.Range("N135").Value = -127390

$ws.Range("H136").Value = 1430.2903
$ws.Range("I136").Value = 1051.826
$ws.Range("K136").Value = 3155.478
$ws.Range("M136").Value = -605.4780000000001
